$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with default (unstyled) General/text formatting, used to restore style
# after temporarily forcing Text format on numeric-looking cells so Excel keeps them as strings.
$fmtSource = $ws.Range("B2")

$ws.Range("D2").Value = '67.780.81'
$ws.Range("E2").Value = '  -1.31%  '

$ws.Range("D3").Value = '3.329.99'
$ws.Range("E3").Value = '  -1.30%  '

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)


$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '582.55'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E5").Value = '  -1.86%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '176.29'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E6").Value = '  -5.42%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  -1.58%  '

$ws.Range("D9").Value = '3.325.82'
$ws.Range("E9").Value = '  -1.18%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.177'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E10").Value = '  -3.01%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '45.64'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E12").Value = '  -3.49%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '661.15'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E14").Value = '  +3.19%  '

$ws.Range("D15").Value = '3.868.49'
$ws.Range("E15").Value = '  -1.05%  '

$ws.Range("E16").Value = '  -1.37%  '

$ws.Range("D17").Value = '67.882.16'
$ws.Range("E17").Value = '  -1.30%  '

$ws.Range("E18").Value = '  -1.09%  '

$ws.Range("D19").Value = '3.333.07'
$ws.Range("E19").Value = '  -1.21%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '17.44'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E20").Value = '  -2.84%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '10.95'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)


$ws.Range("E22").Value = '  -2.31%  '

$ws.Range("E23").Value = '  +7.00%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '17.09'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E24").Value = '  -4.78%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '99.38'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E25").Value = '  -0.04%  '

$ws.Range("E26").Value = '  -5.83%  '

$ws.Range("E27").Value = '  -6.17%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '9.29'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E28").Value = '  -4.78%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '33.53'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E29").Value = '  +1.81%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '7.41'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E30").Value = '  +8.32%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '8.45'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E31").Value = '  -2.69%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '589.98'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E32").Value = '  -4.32%  '

$ws.Range("E33").Value = '  -1.26%  '

$ws.Range("E34").Value = '  -1.19%  '

$ws.Range("D35").Value = '3.731.90'
$ws.Range("E35").Value = '  -6.41%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E36").Value = '  +0.08%  '

$ws.Range("E37").Value = '  +1.43%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '3.35'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E38").Value = '  -9.48%  '

$ws.Range("E39").Value = '  +0.27%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '33.70'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E40").Value = '  +0.01%  '

$ws.Range("E41").Value = '  -5.60%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '3.12'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E42").Value = '  -5.91%  '

$ws.Range("E43").Value = '  -2.72%  '

$ws.Range("E44").Value = '  -5.45%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '3.25'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E45").Value = '  -5.04%  '

$ws.Range("E46").Value = '  -3.85%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '2.60'
$fmtSource.Copy()
$cell.PasteSpecial(-4122)

$ws.Range("E47").Value = '  -0.05%  '

$ws.Range("E48").Value = '  -1.60%  '

$ws.Range("E49").Value = '  +0.06%  '

$ws.Range("E50").Value = '  -0.39%  '

$ws.Range("E51").Value = '  -3.20%  '

$excel.CutCopyMode = $false
